$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Time" row: label in A12, time-of-day value (entered the way Excel
# keeps high-precision fractional-day values, e.g. pasted/typed with E
# notation) in B12, formatted with the built-in h:mm:ss time format.
$ws.Range("A12").Value = "Time"
$ws.Range("B12").Value = 0.097337962962962959
$ws.Range("B12").NumberFormat = "h:mm:ss"

# After entering the value and pressing Enter, the active cell moves one
# row down.
[void]$ws.Range("B13").Select()
